$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ECL")

$ws.Range("D16").Value = 0.4242
$ws.Range("E16").Value = 0.4325
$ws.Range("F16").Value = 0.4409
$ws.Range("G16").Value = 0.4391

$ws.Range("D17").Value = 0.1226
$ws.Range("E17").Value = 0.1319
$ws.Range("F17").Value = 0.1484
$ws.Range("G17").Value = 0.1469

$ws.Range("D18").Value = 0.1035
$ws.Range("E18").Value = 0.1208
$ws.Range("F18").Value = 0.139
$ws.Range("G18").Value = 0.1378

$ws.Range("D19").Value = -0.0901
$ws.Range("E19").Value = -0.0706
$ws.Range("F19").Value = 0.1223
$ws.Range("G19").Value = 0.1241

$ws.Range("D20").Value = 0.1209
$ws.Range("E20").Value = 0.1339
$ws.Range("F20").Value = 0.1359
$ws.Range("G20").Value = 0.1351

$ws.Range("D31").Value = 0.2161
$ws.Range("E31").Value = 0.2083
$ws.Range("F31").Value = 0.214
$ws.Range("G31").Value = 0.2086

$ws.Range("D32").Value = 0.17
$ws.Range("E32").Value = 0.1861
$ws.Range("F32").Value = 0.191
$ws.Range("G32").Value = 0.1927
